$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

# ---------------------------------------------------------------------
# New header cells (first batch): H1, I1, J1
# ---------------------------------------------------------------------
$ws.Range("H1").Value = "id mau sac"
$ws.Range("I1").Value = "id size"
$ws.Range("J1").Value = "số lượng size"

# ---------------------------------------------------------------------
# Row 2 / Row 3 numeric columns A-G
# ---------------------------------------------------------------------
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 20000
$ws.Range("F2").Value = 20000
$ws.Range("G2").Value = 54

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = 121
$ws.Range("F3").Value = 213
$ws.Range("G3").Value = 2

# ---------------------------------------------------------------------
# H2 / I2 (id mau sac, id size lists)
# ---------------------------------------------------------------------
$ws.Range("H2").Value = "1,2,3"
$ws.Range("I2").Value = "1,2,3,4"

# ---------------------------------------------------------------------
# H3 (id mau sac list for row 3)
# ---------------------------------------------------------------------
$ws.Range("H3").Value = "1,2,4"

# ---------------------------------------------------------------------
# J2 / J3 (so luong size)
# ---------------------------------------------------------------------
$ws.Range("J2").Value = 3
$ws.Range("J3").Value = 3

# ---------------------------------------------------------------------
# L1 header ("anh mau sac")
# ---------------------------------------------------------------------
$ws.Range("L1").Value = "ảnh màu sắc"

# ---------------------------------------------------------------------
# L2 / L3 hyperlinks - color-image urls (same text/address for both rows)
# ---------------------------------------------------------------------
$colorImgUrls = "https://royalhelmet.com.vn/ckfinder/userfiles/images/products/Zt7RKI_MG_3349.jpg ,https://royalhelmet.com.vn/ckfinder/userfiles/images/products/P0jXDB_MG_3354.jpg,https://royalhelmet.com.vn/ckfinder/userfiles/images/products/K6uk81_MG_3354-muc.jpg"
$ws.Hyperlinks.Add($ws.Range("L2"), $colorImgUrls) | Out-Null
$ws.Hyperlinks.Add($ws.Range("L3"), $colorImgUrls) | Out-Null

# ---------------------------------------------------------------------
# K1 header ("mo ta mau sac") - added after the L column work
# ---------------------------------------------------------------------
$ws.Range("K1").Value = "mô tả màu sắc"

# ---------------------------------------------------------------------
# I3 ("1,2,4,4") - filled in later
# ---------------------------------------------------------------------
$ws.Range("I3").Value = "1,2,4,4"

# ---------------------------------------------------------------------
# M1 header ("anh chinh")
# ---------------------------------------------------------------------
$ws.Range("M1").Value = "ảnh chính"

# ---------------------------------------------------------------------
# M2 hyperlink - display caches only the first url, cell text is the full list
# ---------------------------------------------------------------------
$m2Display = "https://royalhelmet.com.vn/ckfinder/userfiles/images/products/klIZba_MG_3349.jpg"
$ws.Hyperlinks.Add($ws.Range("M2"), $m2Display, "", "", $m2Display) | Out-Null
$ws.Range("M2").Value = "https://royalhelmet.com.vn/ckfinder/userfiles/images/products/klIZba_MG_3349.jpg, https://royalhelmet.com.vn/ckfinder/userfiles/images/products/TLXYEi_MG_3350.jpg,https://royalhelmet.com.vn/ckfinder/userfiles/images/products/NOYMt4_MG_3352.jpg"

# ---------------------------------------------------------------------
# M3 hyperlink - display caches the whole (matching) url list
# ---------------------------------------------------------------------
$m3Text = "https://royalhelmet.com.vn/ckfinder/userfiles/images/products/b4hEP5SwGTU0royal-m139-v7.jpg,https://royalhelmet.com.vn/ckfinder/userfiles/images/products/3IqhWcroyal-m139-v7-1.jpg,https://royalhelmet.com.vn/ckfinder/userfiles/images/products/hD1uIXroyal-m139-v7-2.jpg"
$ws.Hyperlinks.Add($ws.Range("M3"), $m3Text, "", "", $m3Text) | Out-Null

# ---------------------------------------------------------------------
# Column widths for the new columns (closest achievable values)
# ---------------------------------------------------------------------
$ws.Columns.Item(8).ColumnWidth = 10.857142857142858
$ws.Columns.Item(10).ColumnWidth = 13.571428571428571
$ws.Columns.Item(11).ColumnWidth = 13.571428571428571
$ws.Columns.Item(12).ColumnWidth = 12.857142857142858

# ---------------------------------------------------------------------
# Page setup orientation
# ---------------------------------------------------------------------
$ws.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# Selection / active cell
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("K17").Select() | Out-Null
